$wb = $excel.ActiveWorkbook

# --- YDS and ST sheets: append Week 15 simulation data to long space-separated number strings ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "2 -3 8 0 2 5 7 1 20 -3 4 -2 0 1 0 1 0 3 -1 1 5 3 1 3 1 0 3 1 2 1 0 -1 6 9 8 9 -1 4 8 4 7 3 8 0 13 11 2 2 9 -1 11 1 5 -5 3 22 3 3 2 0 3 3 0 9 2 6 2 7 6 6 13 3 0 5 4 3 1 1 10 2 2 -1 3 7 -1 3 9 1 -1 2 9 5 3 7 3 2 0 1 1 5 0 3 7 -2 3 0 1 4 4 4 6 2 4 1 14 5 0 0 0 6 0 0 3 8 -1 2 4 11 6 2 6 1 4 2 3 1 5 5 8 1 40 2 2 7 12 6 3 5 13 5 9 -1 4 5 8 5 12 5 3 1 5 3 10 2 1 3 3 3 4 -3 11 2 2 5 1 10 -1 6 11 1 1 5 4 3 7 10 3 6 27 2 2 -3 1 2 3 1 12 3 2 4 2 2 5 1 5 14 6 2 6 3 1 3 6 5 4 1 19 16 3 2 1 2 6 2 4 16 10 2 3 13 8 -1 -5 1 3 11 4 2 2 6 1 5 5 3 3 3 6 13 6 10 1 2 -1 4 1 3 2 23 1 15 2 7 0 1 3 37 1 3 8 1 13 1 3 0 8 5 4 0 1 1 1 3 2 3 3 2 -2 0 0 3 11 3 5 3 3 7 -2 2 3 9 17 9 0 17 -1 5 3 5 3 1 0 2 2 3 11 -1 4 9 1 2 2 4 1 4 12 -1 1 6 7 2 5 9 3 2 5 7 5 10 6 13 3 4 0 22 9 3 -1 -6 3 4 7 12 3 3 1 6 18 0 -6 -3 11 0 4 3 0 3 1 3 12 2 5 2 13 0 4 9 5 1 3 -4 4 -1 2 10 2 2 5 2 11 4 1 10 3 3 3 3 -1 13 3 -1 8 4 4 4 27 -3 2 3 6 -1 8 3 5 6 2 6 6 6 6 3 0 -2 0 6 4 5 2 5 9 1 6 2 5 0 5 1 10 2 2 2 2 7 14 3 1 3 1 11 3 8 2 4 2 7 -1 4 1 -1 11 4 -4 7 3 1 7 4 7 6 0 1 0 5 7 4 6 5 3 6 2 8 6 6 -2 20 2 5 6 8 10 2 1 5 1 3 10 -1 1 -1 5 8 3 0 1 8 5 1 -7 13 5 3 3 1 4 1 2 3 3 2 5 6 16 1 9 2 11 2 -2 -2 16 2 1 4 8 1 3 5 3 3 2 2 15 1 7 8 6 4 4 11 2 5 3 8 5 4 5 2 2 6 1 4 38 2 6 4 5 3 0 2 6 27 12 5 2 8 7 13 0 2 10 4 3 2 2 5 3 6 5 0 4 4 7 0 4 5 13 12 12 4 4 1 3 1 5 -1 -8 10 5 2 3 15 4 4 -5 -1 -3 12 9 1 1 1 2 3 4 7 3 2 4 2 -1 5 0 3 4 3 10 2 -2 4 -3 15 5 2 11 1 8 3 4 0 1 4 5 11 20 2 0 11 9 4 9 4 2 8 -2 0 14 3 11 0 0 7 6 7 3 8 4 5 1 0 4 3 4 4 3 1 1 0 1 5 7 3 4 4 3 10 10 2 1 6 3 0 0 7 4 1 7 5 17 1 2 0 2 2 22 1 3 1 16 2 6 2 2 4 1 9 5 -4 -3 13 4 1 0 -6 7 2 3 8 1 6 0 4 3 2 0 1 3 3 11 2 4 6 5 0 13 6 5 12 8 4 6 7 2 0"
$ws.Range("C2").Value = "5 3 1 7 2 2 2 1 -3 -5 13 7 4 10 3 2 3 4 13 8 2 -2 7 4 2 14 4 3 -3 6 0 2 3 3 5 1 21 4 6 4 1 4 16 10 4 3 1 2 2 0 1 2 4 8 3 3 2 1 16 1 -1 11 4 4 5 1 5 -7 6 25 0 9 4 2 0 20 5 3 7 0 -1 4 11 -4 4 1 8 1 2 5 1 1 5 1 50 1 2 7 2 3 1 3 4 0 13 8 6 9 11 -2 0 14 5 4 3 1 1 -1 7 2 6 2 8 3 3 3 12 -1 2 0 0 2 4 9 0 0 5 -2 6 4 1 46 -6 0 5 2 14 -3 -3 5 4 49 0 -1 9 5 8 9 -1 1 8 3 4 0 3 7 2 7 1 4 2 6 2 5 1 1 4 5 0 4 0 1 12 1 6 4 6 5 2 0 6 0 0 4 6 0 7 18 5 6 4 2 4 2 5 2 19 11 5 5 8 1 0 12 3 5 5 3 20 1 -2 2 2 6 -2 4 5 11 16 1 9 2 13 4 0 0 3 0 2 2 2 3 -1 17 9 6 1 4 3 5 10 0 1 7 7 0 3 2 8 0 2 19 2 -5 4 8 -1 11 2 2 2 2 6 3 8 3 -2 5 4 3 0 0 7 -1 2 4 0 0 1 1 9 1 -2 9 2 5 3 9 9 4 11 1 1 6 1 0 6 -2 2 -2 6 4 8 1 14 3 0 9 0 1 10 11 5 5 10 5 5 -3 4 3 4 1 6 38 3 13 9 50 0 6 0 3 -4 6 4 0 4 1 7 2 3 7 5 45 1 3 0 3 3 1 7 3 7 0 1 -2 -2 -1 3 1 -3 3 2 6 2 0 0 0 9 1 8 1 6 0 3 6 7 4 3 20 1 4 3 6 3 1 1 5 0 4 9 7 7 -2 6 0 4 3 22 3 8 4 3 2 3 3 6 6 5 17 8 2 3 2 13 1 11 13 3 3 1 7 -3 1 4 3 2 0 7 2 2 0 14 5 2 3 0 4 1 2 9 -3 5 2 15 -2 2 6 41 6 2 7 5 3 0 9 -1 -6 2 -9 7 2 2 5 11 1 1 2 0 4 6 8 1 -2 4 4 8 2 3 4 3 2 2 0 3 0 7 6 5 1 6 9 1 2 5 2 2 8 2 3 2 4 -3 10 2 5 -3 1 2 -3 3 3 13 5 2 0 13 1 2 -3 4 17 2 1 5 3 1 23 2 10 5 3 12 3 3 0 11 2 1 8 3 -2 1 7 2 0 6 3 3 17 3 6 1 6 4 0 5 11 4 3 3 3 10 5 2 -1 1 3 0 3 9 15 2 2 13 2 1 5 3 2 0 2 0 -1 8 4 2 6 2 2 7 1 1 7 2 -1 5 3 0 11 14 1 7 1 -2 5 6 0 3 5 8 2 1 3 4 8 6 8 6 3 0 -3 4 3 8 3 4 9 6 24 2 4 0 3 6 10 -1 2 7 5 0 0 1 3 0 1 12 4 3 7 1 3 0 17 7 3 15 -1 6 2 1 1 10 1 4 3 0 4 0 5 5 4 0 2 2 7 4 5 4 11 1 4 10 10 6 1 3 2 3 3 4 3 4 1 5 5"
$ws.Range("B3").Value = "16 2 14 21 11 7 6 6 1 17 19 4 3 8 14 10 19 4 8 2 25 9 -3 3 5 13 18 3 15 10 7 26 10 24 14 1 30 8 7 12 21 6 17 14 33 -5 -5 5 12 4 26 6 5 11 -1 14 10 9 17 -2 3 11 8 18 11 1 13 -1 40 4 3 12 0 7 4 4 19 4 13 8 8 11 11 9 16 4 39 3 3 12 8 18 10 6 5 9 6 12 9 5 3 -2 -1 4 1 9 4 3 29 10 8 11 20 8 12 2 0 8 7 5 11 6 4 5 10 11 5 18 5 3 13 6 8 22 7 1 9 5 26 52 6 7 15 13 5 22 12 8 6 21 7 11 11 12 6 0 7 6 14 32 45 8 6 33 7 15 3 6 8 8 68 8 5 8 7 10 11 4 22 19 9 7 3 7 0 8 6 4 10 6 8 7 3 13 6 11 1 20 19 14 4 11 -1 13 8 8 23 27 5 26 5 11 11 11 10 0 8 42 2 3 4 11 7 14 3 25 12 11 7 7 6 -1 7 28 4 6 4 6 9 8 5 5 12 5 9 15 9 20 8 10 8 12 8 8 6 0 5 2 12 6 2 16 3 30 11 1 4 5 31 9 -2 30 6 6 6 16 8 15 7 6 29 -2 3 2 11 13 5 9 9 5 13 12 8 -1 9 5 5 -2 0 0 9 1 11 2 11 20 13 3 10 4 18 -2 5 30 5 1 0 9 7 17 17 3 12 13 5 10 9 6 6 5 2 10 9 5 11 5 -2 4 17 12 0 8 50 14 3 28 7 5 3 8 4 19 7 4 5 9 23 7 13 29 6 1 10 18 16 5 5 6 5 4 5 4 15 6 6 6 13 5 11 6 3 6 5 36 24 18 10 5 5 8 12 7 12 29 12 19 6 1 8 23 4 3 5 1 16 11 16 10 2 6 5 12 7 11 9 2 34 11 17 17 -2 4 7 5 22 12 4 4 16 11 8 12 0 24 12 7 6 10 7 4 4 8 18 9 0 3 10 7 56 19 2 7 9 5 7 6 0 11 17 10 73 14 7 4 15 1 11 15 37 4 1 6 9 3 19 -1 21 33 16 20 11 2 -2 11 23 10 17 24 19 4 30 20 9 18 8 7 17 13 7 16 10 7 -7 9 3 6 4 6 9 26 27 9 13 6 3 7 12 2 7 3 5 6 6 6 4 16 39 10 14 4 -4 -4 9 12 4 -2 7 6 12 40 0 5 5 14 4 17 28 9 12 5 7 10 3 11 2 19 10 3 12 13 6 15 8 27 17 2 13 -2 20 4 12 1 17 14 7 5 14 4 23 20 19 5 22 14 3 1 8 14 13 4 8 11 1 20 22 7 9 8 14 1 15 6 18 21 9 6 16 6 5 6 6 30 7 39 6 4 5 12 12 7 13 4 15 10 18 5 4 15 2 8 6 0 1 1 12 9 8 14 9 16 10 20 2 9 11 11 6 3 12 6 8 23 -1 13 6 7 4 7 6 5 6 5 6 0 10 35 10 3 4 10 9 10 12 6 14 5 12 2 3 19 9 43 -1 6 10 8 19 19 7"
$ws.Range("C3").Value = "16 19 12 5 7 55 -2 0 11 7 10 34 -4 12 23 6 3 27 6 -6 10 4 7 8 4 9 12 4 2 10 8 5 54 3 10 8 5 49 0 6 3 -1 9 4 25 6 13 8 20 10 12 7 6 18 2 10 11 8 9 15 13 3 15 9 10 1 33 31 3 25 13 9 18 22 4 4 10 10 1 17 17 -6 9 23 7 4 56 12 6 8 40 3 5 10 6 49 2 29 11 21 8 10 23 5 9 -4 1 9 18 1 11 5 4 14 10 4 18 6 6 8 32 7 6 6 3 50 4 2 7 8 21 8 10 16 8 7 3 7 12 4 8 7 9 6 0 5 55 4 9 27 19 8 6 2 4 8 26 4 5 18 15 5 0 13 9 5 10 10 9 18 9 5 6 6 5 13 -5 11 12 12 19 8 9 25 5 15 11 2 4 5 8 5 11 14 5 25 3 10 7 54 7 0 9 7 9 9 1 4 1 19 3 7 7 7 2 6 4 6 6 3 0 30 13 4 5 -2 17 4 15 8 8 3 -4 30 50 6 2 1 1 14 6 16 8 5 7 7 2 5 4 9 22 16 6 6 5 17 18 4 11 8 12 12 16 10 -1 13 5 16 6 9 6 12 14 7 19 6 3 5 1 15 7 6 6 9 3 10 3 3 4 15 11 7 7 9 8 1 16 -1 1 44 0 14 13 10 12 14 14 10 4 22 3 3 15 9 7 4 15 30 7 0 5 10 5 6 6 15 17 12 15 36 4 5 27 22 19 24 13 13 13 7 23 17 8 20 19 35 16 10 6 9 11 -3 7 2 11 4 18 16 10 22 1 6 17 4 17 11 7 11 23 6 19 17 3 17 7 19 20 9 12 16 13 4 4 8 2 6 17 3 4 2 8 12 19 33 13 21 10 8 21 13 23 11 13 1 28 3 7 9 7 17 4 7 14 18 2 41 7 13 22 11 5 11 9 5 6 13 10 11 6 5 19 4 13 15 7 15 10 5 12 42 8 8 17 6 12 15 5 7 14 14 6 7 8 7 5 21 4 11 72 20 11 10 12 49 14 -4 11 12 6 32 19 6 20 9 19 27 49 3 5 7 27 19 17 14 9 10 5 13 9 5 12 13 2 9 0 3 11 5 21 19 6 12 6 24 15 3 4 12 5 0 12 1 10 3 15 11 7 13 6 19 5 10 2 20 4 10 26 4 25 17 1 31 5 31 12 19 19 15 7 3 4 -4 9 2 11 14 9 18 7 2 3 29 9 8 3 6 7 11 4 5 14 15 8 6 2 15 4 0 4 3 40 8 13 10 5 3 11 2 19 7 9 4 10 16 27 7 6 4 8 11 7 2 55 6 13 9 39 13 8 2 7 6 1 12 1 4 13 13 4 7 32 1 8 9 3 10 5 3 9 7 5 0 18 10 -5 5 2 15 9 34 7 3 28 14 13 16 2 11 7 4 18 7 15 3 12 24 7 4 0 12 13 9 14 20 12 5 11 5 10 2 4"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("D3").Value = "47 51 48 43 49 56 51 51 53 52 30 60 38 32 55 60 45 51 49 41 43 47 41 53 42 61 52 44 58 57 56 55 38 38 54 56 30 47 52 55 55 55 35 60 57 45 28 60 36 52 54 58 55 48 31 54 54 52 46 44 54 41 40 49 65 45 45 47 48 48 28 50 28 52 53 43 56 45 37 55 48 50 47 52 42 55 38 48 39 45 56 51 47 66 60 55 51 66 44 45 56 49 34 48 37 53 47 33 54 47 43 53 59 57 55 53 42 41 49"
$ws.Range("B4").Value = "68 77 61 66 70 56 59 65 64 50 67 63 61 41 66 65 70 68 73 58 54 66 61 64 62 62 60 64 57 58"
$ws.Range("D4").Value = "0 0 8 0 0 0 11 4 1 0 0 0 0 0 11 0 0 10 5 13 0 0 0 16 3 0 8 2 10 0 3 0 0 0 6 0 0 6 0 10 9 0 0 9 0 0 0 17 0 11 0 0 2 0 0 3 0 10 0 15 3 0 0 0 0 15 3 3 0 8 0 7 0 9 0 0 0 3 0 11 8 0 0 6 33 0 0 5 2 4 13 14 2 0 31 25 0 15 10 4 10 0 0 0 0 0 0 0 6 0 0 0 0 4 0 1 0 0 0"
$ws.Range("B5").Value = "20 31 12 67 30 26 32 29 18 8 19 23 20 0 24 32 22 27 23 17 19 34 21 27 27 27 19 24 21 15"
$ws.Range("D5").Value = "2 0 0 11 4 12 4 1 0 0 1 0 19 3 -3 0 0 13 0 0 0 0 2 0 0 0 3 0 0 0 0 12 0 0 0 0 10 0 3 0 0 12 10 0 14 0 0 0 0 0 22 0 0 0 14 0 0 0 0 0 11 -1 0 0 0 0 0 20 0 0 0 0 0 14 8 0 0 0 0 0 0 6 0 0 0 0 0 10 0 0 0 0 0 0 12 0 0 10 0 0 0 0 0 0 0 0 0 0 0 8 0 0 16 0 13 0 0 0 0 0 0"
$ws.Range("B6").Value = "14 13 10 20 28 23 22 23 4 21 14 25 21 24 20 19 46 24 17 21 25 18 23 29 15 28 22 20 22 22 31 16 31 18 21 27 16 23 21 0 19 28 23 1 22 26 15 23 19 30 24 25 23 28 19 23 30 12 31 28 20"

# --- Numeric summary cells updated on OFF, DEF, ST, TURNS, PEN sheets ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 342
$ws.Range("F2").Value = 134
$ws.Range("G2").Value = 101
$ws.Range("J2").Value = 69
$ws.Range("L2").Value = 592
$ws.Range("M2").Value = 381
$ws.Range("O2").Value = 46
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 1045
$ws.Range("C3").Value = 374
$ws.Range("E3").Value = 59
$ws.Range("F3").Value = 189
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 123
$ws.Range("J3").Value = 101
$ws.Range("N3").Value = 44

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 337
$ws.Range("F2").Value = 116
$ws.Range("G2").Value = 94
$ws.Range("H2").Value = 16
$ws.Range("J2").Value = 57
$ws.Range("L2").Value = 551
$ws.Range("M2").Value = 344
$ws.Range("Q2").Value = 997
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 341
$ws.Range("E3").Value = 63
$ws.Range("F3").Value = 193
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 99
$ws.Range("J3").Value = 128
$ws.Range("N3").Value = 49

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 136
$ws.Range("D2").Value = 119
$ws.Range("F2").Value = 94
$ws.Range("H2").Value = 13
$ws.Range("I2").Value = 2
$ws.Range("B3").Value = 106

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 17
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 16

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 27
$ws.Range("D2").Value = 23
$ws.Range("B3").Value = 26
$ws.Range("B4").Value = 3
$ws.Range("D4").Value = 13
